$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.469.14'
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.790.05'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.30'
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5872'
$ws.Range("E6").Value = '  -1.44%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2760'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.23'
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06719'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07553'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.793.42'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.774'
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6083'
$ws.Range("E14").Value = '  -2.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.032.27'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '75.38'
$ws.Range("E16").Value = '  -3.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008773'
$ws.Range("E17").Value = '  -8.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.455.05'
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.415'
$ws.Range("E19").Value = '  -4.60%  '
$ws.Range("E20").Value = '  -0.11%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '208.37'
$ws.Range("E21").Value = '  -5.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.40'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.785'
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.19'
$ws.Range("E25").Value = '  -2.12%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.971'
$ws.Range("E26").Value = '  +0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1255'
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.35'
$ws.Range("E28").Value = '  -0.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.412'
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06139'
$ws.Range("E30").Value = '  -6.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.418'
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.773'
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.753'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.705'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.046'
$ws.Range("E35").Value = '  -3.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6414'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.503'
$ws.Range("E38").Value = '  -1.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.148.39'
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.309'
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01673'
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("E42").Value = '  -2.45%  '
$ws.Range("E43").Value = '  -0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.55'
$ws.Range("E44").Value = '  +0.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.940.81'
$ws.Range("E45").Value = '  -1.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.98'
$ws.Range("E46").Value = '  -2.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000110'
$ws.Range("E47").Value = '  -4.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.575'
$ws.Range("E48").Value = '  +0.42%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.385'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05430'
$ws.Range("E50").Value = '  -1.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4471'
$ws.Range("E51").Value = '  -1.89%  '
